$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values would otherwise be
# auto-detected as numbers by Excel, so they stay stored as text
# (matching the source data, which is all plain/display-formatted text).
$textCells = @("D5", "D6", "D8", "D9", "D10", "D14", "D15", "D16", "D18", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D29", "D34", "D36", "D37", "D38", "D41", "D42", "D43", "D44", "D45", "D47", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values from the latest crypto-price scrape.
$ws.Range("D2").Value = "35.146.18"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "1.813.96"
$ws.Range("E3").Value = "  -1.67%  "
$ws.Range("E4").Value = "  +0.77%  "
$ws.Range("D5").Value = "232.60"
$ws.Range("E5").Value = "  +2.22%  "
$ws.Range("D6").Value = "0.612"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("D8").Value = "40.89"
$ws.Range("E8").Value = "  -3.44%  "
$ws.Range("D9").Value = "0.322"
$ws.Range("E9").Value = "  +5.51%  "
$ws.Range("D10").Value = "0.0683"
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "2.077.44"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").Value = "1.826.88"
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("D14").Value = "11.04"
$ws.Range("E14").Value = "  -4.01%  "
$ws.Range("D15").Value = "0.658"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").Value = "4.65"
$ws.Range("E16").Value = "  -2.16%  "
$ws.Range("D17").Value = "35.087.99"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "69.50"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "0.0₃0790"
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").Value = "238.70"
$ws.Range("E20").Value = "  -2.74%  "
$ws.Range("D21").Value = "11.89"
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("D22").Value = "4.66"
$ws.Range("E22").Value = "  -0.99%  "
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("D24").Value = "2.26"
$ws.Range("E24").Value = "  +3.81%  "
$ws.Range("D25").Value = "172.71"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").Value = "7.82"
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("D27").Value = "17.46"
$ws.Range("E27").Value = "  -2.22%  "
$ws.Range("E28").Value = "  -1.27%  "
$ws.Range("D29").Value = "1.59"
$ws.Range("E29").Value = "  +19.09%  "
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("E31").Value = "  +3.00%  "
$ws.Range("D32").Value = "3.331.00"
$ws.Range("E32").Value = "  +37.10%  "
$ws.Range("E33").Value = "  +3.16%  "
$ws.Range("D34").Value = "3.99"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("E35").Value = "  -7.05%  "
$ws.Range("D36").Value = "1.15"
$ws.Range("E36").Value = "  +5.84%  "
$ws.Range("D37").Value = "92.64"
$ws.Range("E37").Value = "  +2.72%  "
$ws.Range("D38").Value = "0.678"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("E39").Value = "  +0.48%  "
$ws.Range("D40").Value = "1.308.43"
$ws.Range("E40").Value = "  -2.42%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "1.28"
$ws.Range("E41").Value = "  +2.03%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -2.50%  "
$ws.Range("B43").Value = "HuobiToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D43").Value = "2.47"
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "14.53"
$ws.Range("E44").Value = "  -2.76%  "
$ws.Range("D45").Value = "2.30"
$ws.Range("E45").Value = "  -5.27%  "
$ws.Range("E46").Value = "  -2.72%  "
$ws.Range("D47").Value = "6.34"
$ws.Range("E47").Value = "  +4.80%  "
$ws.Range("E48").Value = "  -1.03%  "
$ws.Range("D49").Value = "1.991.94"
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("E50").Value = "  +0.84%  "
$ws.Range("D51").Value = "0.0650"
$ws.Range("E51").Value = "  +5.04%  "
